# Add two new columns "I0" (I) and "IF" (J) to the sheet, mirroring the
# header styling used by the existing header row and filling in the
# per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they match the rest of the header row (bold, bordered,
# centered).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# --- Data ----------------------------------------------------------------
# Values for column I ("I0") and column J ("IF") for rows 2 through 56.
$iVals = @(1,6,6,7,5,6,5,8,7,8,7,6,7,8,7,8,7,7,9,8,9,6,9,7,7,8,7,7,9,7,7,7,7,7,4,6,7,6,5,8,8,7,5,6,8,8,2,7,4,8,6,7,8,9,6)
$jVals = @(1,6,6,7,5,6,5,8,7,8,7,6,7,8,7,8,7,7,9,8,9,7,9,7,7,8,7,7,9,7,8,7,7,8,4,6,8,6,6,8,8,7,6,6,8,8,3,7,6,8,6,7,8,9,6)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
